# edit.ps1
# Applies the weekly odds-sheet update for Jogos_da_Semana_FlashScore_2024-11-16.xlsx
# - Updates a handful of odds in existing rows 4 and 15
# - Appends a new match row (row 18) for Rayo Zuliano vs Caracas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: update odds (O4:R4) ---
$ws.Range("O4").Value = 1.18
$ws.Range("P4").Value = 4.5
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 2.3

# --- Row 15: update odds ---
$ws.Range("G15").Value = 2.2
$ws.Range("H15").Value = 3.2
$ws.Range("I15").Value = 3.3
$ws.Range("J15").Value = 2.88
$ws.Range("N15").Value = 10
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 1.44
$ws.Range("T15").Value = 2.63
$ws.Range("U15").Value = 1.8
$ws.Range("V15").Value = 1.91
$ws.Range("W15").Value = 7.5
$ws.Range("AF15").Value = 51
$ws.Range("AG15").Value = 251
$ws.Range("AO15").Value = 12
$ws.Range("AT15").Value = 2.63
$ws.Range("AX15").Value = 19

# --- New row 18: Rayo Zuliano vs Caracas (Venezuela Liga FUTVE) ---
$ws.Range("A18").Value = "rZtf1obm"
$ws.Range("B18").Value = "16/11/2024"
$ws.Range("C18").Value = "16:30"
$ws.Range("D18").Value = "VENEZUELA - LIGA FUTVE"
$ws.Range("E18").Value = "Rayo Zuliano"
$ws.Range("F18").Value = "Caracas"
$ws.Range("G18").Value = 2.5
$ws.Range("H18").Value = 3.1
$ws.Range("I18").Value = 2.72
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 2.07
$ws.Range("L18").Value = 3.25
$ws.Range("M18").Value = 1.03
$ws.Range("N18").Value = 6.65
$ws.Range("O18").Value = 1.34
$ws.Range("P18").Value = 2.75
$ws.Range("Q18").Value = 2
$ws.Range("R18").Value = 1.65
$ws.Range("S18").Value = 1.39
$ws.Range("T18").Value = 2.57
$ws.Range("U18").Value = 1.75
$ws.Range("V18").Value = 1.85
$ws.Range("W18").Value = 7.8
$ws.Range("X18").Value = 12.5
$ws.Range("Y18").Value = 9.5
$ws.Range("Z18").Value = 27
$ws.Range("AA18").Value = 21
$ws.Range("AB18").Value = 32
$ws.Range("AC18").Value = 8.5
$ws.Range("AD18").Value = 6
$ws.Range("AE18").Value = 14.5
$ws.Range("AF18").Value = 70
$ws.Range("AG18").Value = 600
$ws.Range("AH18").Value = 7.8
$ws.Range("AI18").Value = 13
$ws.Range("AJ18").Value = 10.25
$ws.Range("AK18").Value = 32
$ws.Range("AL18").Value = 25
$ws.Range("AM18").Value = 35
$ws.Range("AN18").Value = 4.4
$ws.Range("AO18").Value = 12.5
$ws.Range("AP18").Value = 19.5
$ws.Range("AQ18").Value = 50
$ws.Range("AR18").Value = 80
$ws.Range("AS18").Value = 250
$ws.Range("AT18").Value = 2.57
$ws.Range("AU18").Value = 6.7
$ws.Range("AV18").Value = 55
$ws.Range("AW18").Value = 4.6
$ws.Range("AX18").Value = 14.5
$ws.Range("AY18").Value = 22
$ws.Range("AZ18").Value = 65
$ws.Range("BA18").Value = 100
$ws.Range("BB18").Value = 250
$ws.Range("BC18").Value = 51
$ws.Range("BD18").Value = 51
